$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country ranking reshuffles (labels only; underlying shared-string pool is
#     managed automatically, we just need each row to show the right country) ---

# Nepal overtakes Corea del Sur / Dinamarca / Camerun (rows 64-67), each value
# cascades down one rank and Nepal gets fresh totals on row 64.
$ws.Range("A64").Value = "Nepal"
$ws.Range("A65").Value = "Corea del Sur"
$ws.Range("A66").Value = "Dinamarca"
$ws.Range("A67").Value = "Camerun"

# Madagascar overtakes Sri Lanka (rows 108-109)
$ws.Range("A108").Value = "Madagascar"
$ws.Range("A109").Value = "Sri Lanka"

# Dominica overtakes Fiyi (rows 205-206)
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# Islas Malvinas overtakes Groenlandia (rows 209-210)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"

# --- Updated case numbers ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2596894
$ws.Range("C4").Value = 357
$ws.Range("E4").Value = 1387247
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 128153

# India (row 7)
$ws.Range("B7").Value = 530993
$ws.Range("C7").Value = 1416
$ws.Range("E7").Value = 203868

# Catar (row 23)
$ws.Range("B23").Value = 94413
$ws.Range("C23").Value = 750
$ws.Range("D23").Value = 78702
$ws.Range("E23").Value = 15601

# Belgica (row 28)
$ws.Range("B28").Value = 61295
$ws.Range("C28").Value = 86
$ws.Range("E28").Value = 34622

# Paises Bajos (row 33)
$ws.Range("B33").Value = 50147
$ws.Range("C33").Value = 73

# Kuwait (row 35)
$ws.Range("B35").Value = 44942
$ws.Range("C35").Value = 551
$ws.Range("D35").Value = 35494
$ws.Range("E35").Value = 9100
$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 348

# Nepal (row 64) - fresh totals for the new #1 of the group
$ws.Range("B64").Value = 12772
$ws.Range("C64").Value = 463
$ws.Range("D64").Value = 3013
$ws.Range("E64").Value = 9731
$ws.Range("H64").Value = 28

# Corea del Sur (row 65) - inherits old Nepal-row figures
$ws.Range("B65").Value = 12715
$ws.Range("C65").Value = 62
$ws.Range("D65").Value = 11364
$ws.Range("E65").Value = 1069
$ws.Range("H65").Value = 282

# Dinamarca (row 66)
$ws.Range("B66").Value = 12675
$ws.Range("D66").Value = 11508
$ws.Range("E66").Value = 563
$ws.Range("H66").Value = 604

# Camerun (row 67)
$ws.Range("B67").Value = 12592
$ws.Range("D67").Value = 10100
$ws.Range("E67").Value = 2179
$ws.Range("H67").Value = 313

# Consejo Danes para los Refugiados (row 77)
$ws.Range("B77").Value = 6827
$ws.Range("C77").Value = 137
$ws.Range("D77").Value = 985
$ws.Range("E77").Value = 5685
$ws.Range("G77").Value = 4
$ws.Range("H77").Value = 157

# Madagascar (row 108) - fresh totals for the new leader
$ws.Range("B108").Value = 2078
$ws.Range("C108").Value = 73
$ws.Range("D108").Value = 944
$ws.Range("E108").Value = 1116
$ws.Range("G108").Value = 2
$ws.Range("H108").Value = 18

# Sri Lanka (row 109) - inherits old Madagascar-row figures
$ws.Range("B109").Value = 2033
$ws.Range("D109").Value = 1661
$ws.Range("E109").Value = 361
$ws.Range("H109").Value = 11

# Malta (row 150)
$ws.Range("D150").Value = 636
$ws.Range("E150").Value = 25

# --- Footer timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Junio de 2020 a las 14:03"
